$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 41, pushing the existing rows 41-52 down to 43-54.
$ws.Range("A41:A42").EntireRow.Insert()

# New row 41 data
$ws.Range("A41").Value = 5
$ws.Range("B41").Value = "Macroferia Regional de Talca"
$ws.Range("C41").Value = "Maule"
$ws.Range("D41").Value = 44736
$ws.Range("E41").Value = 7
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100107
$ws.Range("H41").Value = "Otros"
$ws.Range("I41").Value = 100107001
$ws.Range("J41").Value = "Caqui"
$ws.Range("K41").Value = "Mankaki"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 120
$ws.Range("N41").Value = 15000
$ws.Range("O41").Value = 15000
$ws.Range("P41").Value = 15000
$ws.Range("Q41").Value = "$/caja 18 kilos granel"
$ws.Range("R41").Value = "Región del Maule"
$ws.Range("S41").Value = 833
$ws.Range("T41").Value = 18

# New row 42 data
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44736
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107001
$ws.Range("J42").Value = "Caqui"
$ws.Range("K42").Value = "Mankaki"
$ws.Range("L42").Value = "Segunda"
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 13000
$ws.Range("O42").Value = 13000
$ws.Range("P42").Value = 13000
$ws.Range("Q42").Value = "$/caja 18 kilos granel"
$ws.Range("R42").Value = "Región del Maule"
$ws.Range("S42").Value = 722
$ws.Range("T42").Value = 18
